{"js": "// Replace the three-digit \u00f7 one-digit division problems (and their\n// answers) throughout the document's table cells with the new values\n// from the commit. Each old value is unique in the document, so a\n// plain text search-and-replace per pair is safe and order-independent.\nconst replacements = [\n  [\"151\u00f79=16, 7\", \"905\u00f79=100, 5\"],\n  [\"969\u00f79=107, 6\", \"830\u00f79=92, 2\"],\n  [\"301\u00f72=150, 1\", \"319\u00f74=79, 3\"],\n  [\"497\u00f78=62, 1\", \"355\u00f75=71, 0\"],\n  [\"186\u00f79=20, 6\", \"894\u00f76=149, 0\"],\n  [\"393\u00f78=49, 1\", \"505\u00f75=101, 0\"],\n  [\"618\u00f74=154, 2\", \"862\u00f75=172, 2\"],\n  [\"356\u00f75=71, 1\", \"660\u00f73=220, 0\"],\n  [\"478\u00f76=79, 4\", \"428\u00f79=47, 5\"],\n  [\"882\u00f78=110, 2\", \"360\u00f73=120, 0\"],\n  [\"651\u00f72=325, 1\", \"699\u00f78=87, 3\"],\n  [\"458\u00f75=91, 3\", \"633\u00f72=316, 1\"],\n  [\"810\u00f77=115, 5\", \"491\u00f78=61, 3\"],\n  [\"817\u00f75=163, 2\", \"815\u00f75=163, 0\"],\n  [\"402\u00f79=44, 6\", \"993\u00f77=141, 6\"],\n  [\"166\u00f79=18, 4\", \"245\u00f72=122, 1\"],\n  [\"804\u00f78=100, 4\", \"814\u00f79=90, 4\"],\n  [\"926\u00f77=132, 2\", \"406\u00f75=81, 1\"],\n  [\"588\u00f74=147, 0\", \"642\u00f76=107, 0\"],\n  [\"756\u00f72=378, 0\", \"697\u00f74=174, 1\"],\n  [\"215\u00f73=71, 2\", \"787\u00f75=157, 2\"],\n  [\"939\u00f79=104, 3\", \"819\u00f74=204, 3\"],\n  [\"467\u00f78=58, 3\", \"897\u00f73=299, 0\"],\n  [\"887\u00f78=110, 7\", \"894\u00f73=298, 0\"],\n  [\"692\u00f73=230, 2\", \"254\u00f74=63, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit \u00f7 one-digit division problems (and their\n# answers) throughout the document's table cells with the new values\n# from the commit. Each old value is unique in the document, so a\n# plain Find/Replace per pair (wdReplaceAll) is safe and order-independent.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"151\u00f79=16, 7\", \"905\u00f79=100, 5\"),\n    @(\"969\u00f79=107, 6\", \"830\u00f79=92, 2\"),\n    @(\"301\u00f72=150, 1\", \"319\u00f74=79, 3\"),\n    @(\"497\u00f78=62, 1\", \"355\u00f75=71, 0\"),\n    @(\"186\u00f79=20, 6\", \"894\u00f76=149, 0\"),\n    @(\"393\u00f78=49, 1\", \"505\u00f75=101, 0\"),\n    @(\"618\u00f74=154, 2\", \"862\u00f75=172, 2\"),\n    @(\"356\u00f75=71, 1\", \"660\u00f73=220, 0\"),\n    @(\"478\u00f76=79, 4\", \"428\u00f79=47, 5\"),\n    @(\"882\u00f78=110, 2\", \"360\u00f73=120, 0\"),\n    @(\"651\u00f72=325, 1\", \"699\u00f78=87, 3\"),\n    @(\"458\u00f75=91, 3\", \"633\u00f72=316, 1\"),\n    @(\"810\u00f77=115, 5\", \"491\u00f78=61, 3\"),\n    @(\"817\u00f75=163, 2\", \"815\u00f75=163, 0\"),\n    @(\"402\u00f79=44, 6\", \"993\u00f77=141, 6\"),\n    @(\"166\u00f79=18, 4\", \"245\u00f72=122, 1\"),\n    @(\"804\u00f78=100, 4\", \"814\u00f79=90, 4\"),\n    @(\"926\u00f77=132, 2\", \"406\u00f75=81, 1\"),\n    @(\"588\u00f74=147, 0\", \"642\u00f76=107, 0\"),\n    @(\"756\u00f72=378, 0\", \"697\u00f74=174, 1\"),\n    @(\"215\u00f73=71, 2\", \"787\u00f75=157, 2\"),\n    @(\"939\u00f79=104, 3\", \"819\u00f74=204, 3\"),\n    @(\"467\u00f78=58, 3\", \"897\u00f73=299, 0\"),\n    @(\"887\u00f78=110, 7\", \"894\u00f73=298, 0\"),\n    @(\"692\u00f73=230, 2\", \"254\u00f74=63, 2\")\n)\n\n# wdReplaceAll = 2, wdFindContinue = 1\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll) | Out-Null\n}\n"}
